# Auto-generated edit script: Add violent-crime data for 2024-12-29
# Updates column K (year 2024) totals across Citywide Totals, By Neighborhood,
# and per-neighborhood detail sheets to reflect one additional day of data.

$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 7849  # Aggravated Assault: 7836 -> 7849
$ws.Range('K3').Value = 8136  # Aggravated Battery: 8117 -> 8136
$ws.Range('I4').Value = 1818  # Criminal Sexual Assault: 1817 -> 1818
$ws.Range('K4').Value = 1714  # Criminal Sexual Assault: 1708 -> 1714
$ws.Range('K6').Value = 9055  # Robbery: 9030 -> 9055
$ws.Range('I7').Value = 26278  # Total: 26277 -> 26278
$ws.Range('K7').Value = 27332  # Total: 27269 -> 27332

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 237  # Albany Park: 236 -> 237
$ws.Range('K7').Value = 812  # Auburn Gresham: 810 -> 812
$ws.Range('K8').Value = 1788  # Austin: 1785 -> 1788
$ws.Range('K9').Value = 131  # Avalon Park: 130 -> 131
$ws.Range('K11').Value = 479  # Belmont Cragin: 478 -> 479
$ws.Range('K16').Value = 66  # Bucktown: 65 -> 66
$ws.Range('K19').Value = 784  # Chatham: 783 -> 784
$ws.Range('K22').Value = 85  # Clearing: 84 -> 85
$ws.Range('K23').Value = 267  # Douglas: 265 -> 267
$ws.Range('K29').Value = 1506  # Englewood: 1505 -> 1506
$ws.Range('K34').Value = 155  # Garfield Ridge: 154 -> 155
$ws.Range('K37').Value = 900  # Grand Crossing: 898 -> 900
$ws.Range('K41').Value = 180  # Hermosa: 179 -> 180
$ws.Range('K42').Value = 1016  # Humboldt Park: 1014 -> 1016
$ws.Range('K48').Value = 347  # Lake View: 343 -> 347
$ws.Range('K50').Value = 122  # Lincoln Square: 121 -> 122
$ws.Range('K51').Value = 354  # Little Italy, UIC: 352 -> 354
$ws.Range('K53').Value = 347  # Logan Square: 346 -> 347
$ws.Range('K54').Value = 533  # Loop: 530 -> 533
$ws.Range('K55').Value = 298  # Lower West Side: 297 -> 298
$ws.Range('I63').Value = 240  # NO NEIGHBORHOOD DATA: 239 -> 240
$ws.Range('K63').Value = 80  # NO NEIGHBORHOOD DATA: 78 -> 80
$ws.Range('K65').Value = 634  # New City: 632 -> 634
$ws.Range('K67').Value = 1065  # North Lawndale: 1064 -> 1065
$ws.Range('K70').Value = 50  # O'Hare: 49 -> 50
$ws.Range('K75').Value = 90  # Pullman: 89 -> 90
$ws.Range('K76').Value = 376  # River North: 375 -> 376
$ws.Range('K77').Value = 180  # Riverdale: 179 -> 180
$ws.Range('K80').Value = 104  # Rush & Division: 103 -> 104
$ws.Range('K83').Value = 582  # South Chicago: 580 -> 582
$ws.Range('K84').Value = 223  # South Deering: 221 -> 223
$ws.Range('K85').Value = 1266  # South Shore: 1260 -> 1266
$ws.Range('K90').Value = 260  # Washington Heights: 259 -> 260
$ws.Range('K93').Value = 111  # West Lawn: 110 -> 111
$ws.Range('K94').Value = 365  # West Loop: 364 -> 365
$ws.Range('K95').Value = 454  # West Pullman: 453 -> 454
$ws.Range('K96').Value = 296  # West Ridge: 295 -> 296
$ws.Range('K97').Value = 225  # West Town: 222 -> 225
$ws.Range('K98').Value = 147  # Wicker Park: 146 -> 147
$ws.Range('K99').Value = 464  # Woodlawn: 462 -> 464
$ws.Range('I101').Value = 26278  # Total: 26277 -> 26278
$ws.Range('K101').Value = 27332  # Total: 27269 -> 27332

# --- West Ridge ---
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K6').Value = 119  # Robbery: 118 -> 119
$ws.Range('K7').Value = 296  # Total: 295 -> 296

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K3').Value = 257  # Aggravated Battery: 255 -> 257
$ws.Range('K7').Value = 812  # Total: 810 -> 812

# --- Belmont Cragin ---
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K4').Value = 28  # Criminal Sexual Assault: 27 -> 28
$ws.Range('K7').Value = 479  # Total: 478 -> 479

# --- South Shore ---
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 421  # Aggravated Assault: 418 -> 421
$ws.Range('K3').Value = 437  # Aggravated Battery: 436 -> 437
$ws.Range('K6').Value = 311  # Robbery: 309 -> 311
$ws.Range('K7').Value = 1266  # Total: 1260 -> 1266

# --- Logan Square ---
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K6').Value = 145  # Robbery: 144 -> 145
$ws.Range('K7').Value = 347  # Total: 346 -> 347

# --- Austin ---
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 542  # Aggravated Battery: 541 -> 542
$ws.Range('K6').Value = 602  # Robbery: 600 -> 602
$ws.Range('K7').Value = 1788  # Total: 1785 -> 1788

# --- South Chicago ---
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 202  # Aggravated Assault: 200 -> 202
$ws.Range('K7').Value = 582  # Total: 580 -> 582

# --- West Pullman ---
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 153  # Aggravated Assault: 152 -> 153
$ws.Range('K7').Value = 454  # Total: 453 -> 454

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 258  # Aggravated Assault: 257 -> 258
$ws.Range('K6').Value = 270  # Robbery: 269 -> 270
$ws.Range('K7').Value = 900  # Total: 898 -> 900

# --- New City ---
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 206  # Aggravated Assault: 205 -> 206
$ws.Range('K3').Value = 150  # Aggravated Battery: 149 -> 150
$ws.Range('K7').Value = 634  # Total: 632 -> 634

# --- Woodlawn ---
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 124  # Aggravated Assault: 123 -> 124
$ws.Range('K6').Value = 112  # Robbery: 111 -> 112
$ws.Range('K7').Value = 464  # Total: 462 -> 464

# --- North Lawndale ---
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 292  # Aggravated Assault: 291 -> 292
$ws.Range('K3').Value = 387  # Aggravated Battery: 386 -> 387
$ws.Range('K4').Value = 61  # Criminal Sexual Assault: 62 -> 61
$ws.Range('K7').Value = 1065  # Total: 1064 -> 1065

# --- South Deering ---
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 90  # Aggravated Battery: 89 -> 90
$ws.Range('K6').Value = 49  # Robbery: 48 -> 49
$ws.Range('K7').Value = 223  # Total: 221 -> 223

# --- Loop ---
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K4').Value = 39  # Criminal Sexual Assault: 38 -> 39
$ws.Range('K6').Value = 287  # Robbery: 285 -> 287
$ws.Range('K7').Value = 533  # Total: 530 -> 533

# --- Englewood ---
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K3').Value = 532  # Aggravated Battery: 531 -> 532
$ws.Range('K7').Value = 1506  # Total: 1505 -> 1506

# --- Lake View ---
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K3').Value = 83  # Aggravated Battery: 82 -> 83
$ws.Range('K4').Value = 52  # Criminal Sexual Assault: 51 -> 52
$ws.Range('K6').Value = 158  # Robbery: 156 -> 158
$ws.Range('K7').Value = 347  # Total: 343 -> 347

# --- Chatham ---
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K6').Value = 262  # Robbery: 261 -> 262
$ws.Range('K7').Value = 784  # Total: 783 -> 784

# --- River North ---
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K6').Value = 185  # Robbery: 184 -> 185
$ws.Range('K7').Value = 376  # Total: 375 -> 376

# --- Hermosa ---
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K3').Value = 38  # Aggravated Battery: 37 -> 38
$ws.Range('K7').Value = 180  # Total: 179 -> 180

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 270  # Aggravated Assault: 269 -> 270
$ws.Range('K3').Value = 297  # Aggravated Battery: 295 -> 297
$ws.Range('K6').Value = 388  # Robbery: 389 -> 388
$ws.Range('K7').Value = 1016  # Total: 1014 -> 1016

# --- Lower West Side ---
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K6').Value = 111  # Robbery: 110 -> 111
$ws.Range('K7').Value = 298  # Total: 297 -> 298

# --- Douglas ---
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K3').Value = 93  # Aggravated Battery: 92 -> 93
$ws.Range('K6').Value = 74  # Robbery: 73 -> 74
$ws.Range('K7').Value = 267  # Total: 265 -> 267

# --- West Lawn ---
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K2').Value = 35  # Aggravated Assault: 34 -> 35
$ws.Range('K7').Value = 111  # Total: 110 -> 111

# --- Garfield Ridge ---
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K3').Value = 42  # Aggravated Battery: 41 -> 42
$ws.Range('K7').Value = 155  # Total: 154 -> 155

# --- West Loop ---
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K3').Value = 75  # Aggravated Battery: 76 -> 75
$ws.Range('K6').Value = 171  # Robbery: 169 -> 171
$ws.Range('K7').Value = 365  # Total: 364 -> 365

# --- Wicker Park ---
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K2').Value = 28  # Aggravated Assault: 27 -> 28
$ws.Range('K7').Value = 147  # Total: 146 -> 147

# --- Lincoln Square ---
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K6').Value = 54  # Robbery: 53 -> 54
$ws.Range('K7').Value = 122  # Total: 121 -> 122

# --- Avalon Park ---
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K3').Value = 48  # Aggravated Battery: 47 -> 48
$ws.Range('K7').Value = 131  # Total: 130 -> 131

# --- Albany Park ---
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K3').Value = 70  # Aggravated Battery: 69 -> 70
$ws.Range('K7').Value = 237  # Total: 236 -> 237

# --- West Town ---
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K4').Value = 8  # Criminal Sexual Assault: 7 -> 8
$ws.Range('K6').Value = 126  # Robbery: 124 -> 126
$ws.Range('K7').Value = 225  # Total: 222 -> 225

# --- O'Hare ---
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('K4').Value = 8  # Criminal Sexual Assault: 7 -> 8
$ws.Range('K7').Value = 50  # Total: 49 -> 50

# --- Pullman ---
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('K3').Value = 30  # Aggravated Battery: 29 -> 30
$ws.Range('K7').Value = 90  # Total: 89 -> 90

# --- Washington Heights ---
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K6').Value = 70  # Robbery: 69 -> 70
$ws.Range('K7').Value = 260  # Total: 259 -> 260

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 98  # Aggravated Battery: 97 -> 98
$ws.Range('K6').Value = 115  # Robbery: 114 -> 115
$ws.Range('K7').Value = 354  # Total: 352 -> 354

# --- Clearing ---
$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('K6').Value = 19  # Robbery: 18 -> 19
$ws.Range('K7').Value = 85  # Total: 84 -> 85

# --- Riverdale ---
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K6').Value = 28  # Robbery: 27 -> 28
$ws.Range('K7').Value = 180  # Total: 179 -> 180

# --- Rush & Division ---
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K6').Value = 53  # Robbery: 52 -> 53
$ws.Range('K7').Value = 104  # Total: 103 -> 104

# --- Bucktown ---
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K4').Value = 5  # Criminal Sexual Assault: 4 -> 5
$ws.Range('K7').Value = 66  # Total: 65 -> 66
